$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("GW_PC_AccountCreation")
$ws2 = $wb.Worksheets.Item("GW_PC_GoogleSearch")
$ws3 = $wb.Worksheets.Item("GW_BC_BillingSummaryAPI")

# --- Sheet1: GW_PC_AccountCreation ---

# Header row: Feature -> UserStory
$ws1.Range("B1").Value = "UserStory"

# Q1 header text stays "Account_ID" but is restyled to the teal header look
$ws1.Range("Q1").Value = "Account_ID"
[void]$ws3.Range("J1").Copy()
[void]$ws1.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats

# Row 3 data cells (row2's G2/H2 are formulas referencing G3/H3, so they follow automatically)
$ws1.Range("G3").Value = "Team"
$ws1.Range("H3").Value = "Demo"
$ws1.Range("O3").Value = "ACV Property Insurance"
$ws1.Range("P3").Value = "301-008578 ACV Property Insurance"
$ws1.Range("Q3").Value = "'4025692771"

# Restyle Q2/Q3 to the teal body look
[void]$ws3.Range("J2").Copy()
[void]$ws1.Range("Q2").PasteSpecial(-4122)
[void]$ws1.Range("Q3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Sheet2: GW_PC_GoogleSearch --- (selection only)
$ws2.Activate()
[void]$ws2.Range("A2:G3").Select()

# --- Sheet3: GW_BC_BillingSummaryAPI --- (selection only)
$ws3.Activate()
[void]$ws3.Range("J1:J2").Select()

# --- Sheet1 stays the active/tabSelected sheet, with B1 selected ---
$ws1.Activate()
[void]$ws1.Range("B1").Select()
